$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: 7,-9)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: -8,8)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: 0,-2)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: -5,-5)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: -10,8)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: 3,9)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: -8,-7)"

$ws.Range("A3").Value = "cost: 616.4732081726651"
$ws.Range("A4").Value = "time: 118.29464163453304"
